$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 17 (the second data row) and shift cells up
$ws.Rows("17").Delete()

# Update totals
$ws.Range("E11").Value = 56940
$ws.Range("F13").Value = 1

# Update the remaining data row's "Periodo Mora" value
$ws.Range("E16").Value = "2509"
$ws.Range("E16").HorizontalAlignment = -4108
